$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.08 = 7502.12 pesos`n✅ 7502.12 pesos = 2.07 = 943.66 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 479.999
$wsTasas.Range("O10").Value = 3601.01
$wsTasas.Range("N12").Value = 3630
$wsTasas.Range("O12").Value = 456.6
